$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number need an explicit Text
# number format first, otherwise Excel auto-converts the literal into a
# numeric value (dropping things like trailing zeros, e.g. "1.010" -> 1.01).

# Row 2
$ws.Range("D2").Formula = '26.583.14'
$ws.Range("E2").Formula = '  +0.68%  '

# Row 3
$ws.Range("D3").Formula = '1.819.63'
$ws.Range("E3").Formula = '  +1.43%  '

# Row 4
$ws.Range("E4").Formula = '  +0.39%  '

# Row 5
$ws.Range("E5").Formula = '  +0.23%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '305.23'
$ws.Range("E6").Formula = '  -0.40%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4659'
$ws.Range("E7").Formula = '  +2.25%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3587'
$ws.Range("E8").Formula = '  -0.96%  '

# Row 9
$ws.Range("E9").Formula = '  +0.13%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8971'
$ws.Range("E10").Formula = '  +2.04%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07788'
$ws.Range("E11").Formula = '  -0.46%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '19.31'
$ws.Range("E12").Formula = '  -1.01%  '

# Row 13
$ws.Range("D13").Formula = '1.829.47'
$ws.Range("E13").Formula = '  +1.61%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.246'
$ws.Range("E14").Formula = '  -0.66%  '

# Row 15
$ws.Range("E15").Formula = '  +0.14%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '87.06'
$ws.Range("E16").Formula = '  +2.39%  '

# Row 17
$ws.Range("E17").Formula = '  +0.25%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008548'
$ws.Range("E18").Formula = '  -0.27%  '

# Row 19
$ws.Range("E19").Formula = '  +0.17%  '

# Row 20
$ws.Range("D20").Formula = '26.629.55'
$ws.Range("E20").Formula = '  +0.70%  '

# Row 21
$ws.Range("E21").Formula = '  -0.99%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.002'
$ws.Range("E22").Formula = '  +0.30%  '

# Row 23
$ws.Range("E23").Formula = '  +0.34%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.929'
$ws.Range("E24").Formula = '  -2.29%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.92'
$ws.Range("E25").Formula = '  -0.48%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '17.86'
$ws.Range("E26").Formula = '  -0.28%  '

# Row 27
$ws.Range("E27").Formula = '  -3.70%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '113.53'
$ws.Range("E28").Formula = '  +1.42%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.792'
$ws.Range("E29").Formula = '  -1.47%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08786'
$ws.Range("E30").Formula = '  +1.50%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.117'
$ws.Range("E31").Formula = '  +1.33%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7267'
$ws.Range("E32").Formula = '  +0.52%  '

# Row 33
$ws.Range("B33").Formula = 'Filecoin'
$ws.Range("C33").Formula = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.419'
$ws.Range("E33").Formula = '  -0.66%  '

# Row 34
$ws.Range("B34").Formula = 'RenderToken'
$ws.Range("C34").Formula = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.711'
$ws.Range("E34").Formula = '  +2.33%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.121'
$ws.Range("E35").Formula = '  +1.10%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.074'
$ws.Range("E36").Formula = '  -0.22%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01923'
$ws.Range("E37").Formula = '  -0.97%  '

# Row 38
$ws.Range("E38").Formula = '  +1.63%  '

# Row 39
$ws.Range("E39").Formula = '  -0.51%  '

# Row 40
$ws.Range("B40").Formula = 'TheSandbox'
$ws.Range("C40").Formula = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5023'
$ws.Range("E40").Formula = '  -3.78%  '

# Row 41
$ws.Range("B41").Formula = 'FraxShare'
$ws.Range("C41").Formula = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.802'
$ws.Range("E41").Formula = '  -1.14%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1488'
$ws.Range("E42").Formula = '  -2.60%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '7.927'
$ws.Range("E43").Formula = '  -1.11%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.010'
$ws.Range("E44").Formula = '  +0.36%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4635'
$ws.Range("E45").Formula = '  -1.22%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.951'
$ws.Range("E46").Formula = '  +0.48%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '97.69'
$ws.Range("E47").Formula = '  -2.29%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.551'
$ws.Range("E48").Formula = '  -2.24%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05986'
$ws.Range("E49").Formula = '  +0.31%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '63.64'
$ws.Range("E50").Formula = '  -0.71%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '35.63'
$ws.Range("E51").Formula = '  -1.73%  '
